$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2..29 (CRN label + VOL,B0,CET_100,CET_300,SP2_PERC,SP,SP2,SP3)
# Row index in $data corresponds to data index 0 => worksheet row 2, etc.
$data = @(
    @{A="crn_00000"; B=3508.135745; C=234.672742; D=0.109738; E=3.344617; F=52.34375; G=9; H=268; I=235}
    @{A="crn_00001"; B=3501.974888; C=241.188986; D=0.155371; E=3.357206; F=52.34375; G=11; H=268; I=233}
    @{A="crn_00002"; B=3459.119351; C=260.416593; D=0.163422; E=4.167658; F=47.265625; G=7; H=242; I=263}
    @{A="crn_00003"; B=3502.943533; C=244.732097; D=0.414841; E=4.176454; F=50; G=6; H=256; I=250}
    @{A="crn_00004"; B=3464.268178; C=229.858607; D=0.158334; E=3.497626; F=52.734375; G=13; H=270; I=229}
    @{A="crn_00005"; B=3460.972341; C=247.056466; D=0.141049; E=4.322701; F=52.34375; G=7; H=268; I=237}
    @{A="crn_00006"; B=3503.012729; C=241.661907; D=-0.0027; E=2.933347; F=48.046875; G=11; H=246; I=255}
    @{A="crn_00007"; B=3431.130834; C=232.349184; D=0.583655; E=4.395775; F=54.296875; G=15; H=278; I=219}
    @{A="crn_00008"; B=3514.997317; C=228.282411; D=0.136162; E=4.701634; F=53.125; G=14; H=272; I=226}
    @{A="crn_00009"; B=3489.605802; C=234.887805; D=0.577572; E=5.479533; F=52.34375; G=13; H=268; I=231}
    @{A="crn_00010"; B=3482.777324; C=232.439478; D=0.152886; E=4.346383; F=52.734375; G=11; H=270; I=231}
    @{A="crn_00011"; B=3469.834646; C=255.414832; D=0.27883; E=4.650076; F=48.828125; G=13; H=250; I=249}
    @{A="crn_00012"; B=3484.983534; C=248.273535; D=0.035614; E=2.849169; F=46.875; G=9; H=240; I=263}
    @{A="crn_00013"; B=3455.826781; C=219.778896; D=-0.13634; E=5.223131; F=54.296875; G=11; H=278; I=223}
    @{A="crn_00014"; B=3530.416387; C=230.878644; D=0.159849; E=3.992783; F=51.5625; G=12; H=264; I=236}
    @{A="crn_00015"; B=3499.070024; C=221.343588; D=-0.422883; E=2.309825; F=54.6875; G=10; H=280; I=222}
    @{A="crn_00016"; B=3470.591046; C=249.0827; D=0.655085; E=5.414664; F=49.21875; G=12; H=252; I=248}
    @{A="crn_00017"; B=3458.158803; C=247.665766; D=-0.146489; E=3.072862; F=46.875; G=9; H=240; I=263}
    @{A="crn_00018"; B=3544.276254; C=228.45318; D=0.238154; E=5.708945; F=51.5625; G=8; H=264; I=240}
    @{A="crn_00019"; B=3484.156096; C=248.48814; D=-0.000191; E=3.276276; F=54.6875; G=10; H=280; I=222}
    @{A="crn_00020"; B=3497.618195; C=220.269765; D=0.351229; E=5.502011; F=52.34375; G=9; H=268; I=235}
    @{A="crn_00021"; B=3539.048435; C=236.221391; D=-0.023666; E=2.982932; F=54.6875; G=7; H=280; I=225}
    @{A="crn_00022"; B=3460.83506; C=237.252454; D=0.524092; E=6.613269; F=50.78125; G=12; H=260; I=240}
    @{A="crn_00023"; B=3456.992661; C=249.17126; D=0.305884; E=5.452449; F=50.78125; G=7; H=260; I=245}
    @{A="crn_00024"; B=3499.346609; C=234.436882; D=0.318476; E=5.399811; F=53.90625; G=13; H=276; I=223}
    @{A="crn_00025"; B=3473.617744; C=234.117044; D=1.364223; E=8.262285; F=54.6875; G=13; H=280; I=219}
    @{A="crn_00026"; B=3449.042456; C=248.919302; D=-0.343943; E=2.277313; F=49.21875; G=10; H=252; I=250}
    @{A="crn_00027"; B=3510.837721; C=234.843471; D=0.213731; E=6.076105; F=51.5625; G=12; H=264; I=236}
)

$startRow = 2
$lastExistingRow = 26
$firstNewRow = 27
$lastNewRow = 29

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rec = $data[$i]

    if ($row -gt $lastExistingRow) {
        # New row: write the CRN label and apply the same thin-border style used by the rest of the table
        $ws.Cells.Item($row, 1).Value2 = $rec.A
    }

    $ws.Cells.Item($row, 2).Value2 = $rec.B
    $ws.Cells.Item($row, 3).Value2 = $rec.C
    $ws.Cells.Item($row, 4).Value2 = $rec.D
    $ws.Cells.Item($row, 5).Value2 = $rec.E
    $ws.Cells.Item($row, 6).Value2 = $rec.F
    $ws.Cells.Item($row, 7).Value2 = $rec.G
    $ws.Cells.Item($row, 8).Value2 = $rec.H
    $ws.Cells.Item($row, 9).Value2 = $rec.I
}

# Apply the same cell formatting (thin border all around, like the rest of the table) to the newly added rows
$newRange = $ws.Range("A" + $firstNewRow + ":I" + $lastNewRow)
$newRange.Borders.LineStyle = 1
$newRange.Borders.Weight = 2
